{"js": "// Replace each three-digit-by-one-digit multiplication problem's text\n// with its newly generated counterpart, per the commit's regenerated\n// answer set. Every occurrence in the document is a unique math-fact\n// string inside a table cell run, so a literal (non-wildcard) search\n// + insertText(..., Word.InsertLocation.replace) on the found range\n// swaps the text while leaving the run's formatting (font, size, etc.)\n// untouched.\nconst replacements = [\n  [\"915\u00d74=3660\", \"689\u00d74=2756\"],\n  [\"786\u00d78=6288\", \"707\u00d75=3535\"],\n  [\"452\u00d78=3616\", \"393\u00d72=786\"],\n  [\"705\u00d74=2820\", \"502\u00d75=2510\"],\n  [\"181\u00d74=724\", \"403\u00d74=1612\"],\n  [\"126\u00d72=252\", \"639\u00d76=3834\"],\n  [\"749\u00d78=5992\", \"680\u00d72=1360\"],\n  [\"686\u00d78=5488\", \"923\u00d79=8307\"],\n  [\"763\u00d74=3052\", \"631\u00d75=3155\"],\n  [\"704\u00d78=5632\", \"396\u00d78=3168\"],\n  [\"240\u00d75=1200\", \"742\u00d78=5936\"],\n  [\"978\u00d75=4890\", \"426\u00d73=1278\"],\n  [\"136\u00d73=408\", \"312\u00d74=1248\"],\n  [\"356\u00d73=1068\", \"410\u00d73=1230\"],\n  [\"874\u00d73=2622\", \"855\u00d76=5130\"],\n  [\"972\u00d72=1944\", \"652\u00d78=5216\"],\n  [\"376\u00d79=3384\", \"690\u00d75=3450\"],\n  [\"646\u00d74=2584\", \"475\u00d78=3800\"],\n  [\"173\u00d79=1557\", \"900\u00d73=2700\"],\n  [\"954\u00d78=7632\", \"198\u00d74=792\"],\n  [\"357\u00d79=3213\", \"267\u00d73=801\"],\n  [\"998\u00d76=5988\", \"458\u00d74=1832\"],\n  [\"606\u00d76=3636\", \"793\u00d74=3172\"],\n  [\"696\u00d76=4176\", \"165\u00d76=990\"],\n  [\"115\u00d79=1035\", \"444\u00d72=888\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each three-digit-by-one-digit multiplication problem's text\n# with its newly generated counterpart, per the commit's regenerated\n# answer set. Every occurrence in the document is a unique math-fact\n# string inside a table-cell run, so Find/Replace (wdReplaceAll) on the\n# whole document body swaps the text while leaving the run's\n# formatting (font, size, etc.) untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"915\u00d74=3660\", \"689\u00d74=2756\"),\n    @(\"786\u00d78=6288\", \"707\u00d75=3535\"),\n    @(\"452\u00d78=3616\", \"393\u00d72=786\"),\n    @(\"705\u00d74=2820\", \"502\u00d75=2510\"),\n    @(\"181\u00d74=724\", \"403\u00d74=1612\"),\n    @(\"126\u00d72=252\", \"639\u00d76=3834\"),\n    @(\"749\u00d78=5992\", \"680\u00d72=1360\"),\n    @(\"686\u00d78=5488\", \"923\u00d79=8307\"),\n    @(\"763\u00d74=3052\", \"631\u00d75=3155\"),\n    @(\"704\u00d78=5632\", \"396\u00d78=3168\"),\n    @(\"240\u00d75=1200\", \"742\u00d78=5936\"),\n    @(\"978\u00d75=4890\", \"426\u00d73=1278\"),\n    @(\"136\u00d73=408\", \"312\u00d74=1248\"),\n    @(\"356\u00d73=1068\", \"410\u00d73=1230\"),\n    @(\"874\u00d73=2622\", \"855\u00d76=5130\"),\n    @(\"972\u00d72=1944\", \"652\u00d78=5216\"),\n    @(\"376\u00d79=3384\", \"690\u00d75=3450\"),\n    @(\"646\u00d74=2584\", \"475\u00d78=3800\"),\n    @(\"173\u00d79=1557\", \"900\u00d73=2700\"),\n    @(\"954\u00d78=7632\", \"198\u00d74=792\"),\n    @(\"357\u00d79=3213\", \"267\u00d73=801\"),\n    @(\"998\u00d76=5988\", \"458\u00d74=1832\"),\n    @(\"606\u00d76=3636\", \"793\u00d74=3172\"),\n    @(\"696\u00d76=4176\", \"165\u00d76=990\"),\n    @(\"115\u00d79=1035\", \"444\u00d72=888\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Forward = $true\n    $find.Wrap = 1            # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute(\n        $find.Text,\n        $false,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $find.Replacement.Text,\n        2                      # wdReplaceAll\n    ) | Out-Null\n}\n"}
